$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 40 - new data row appended to the sheet
$ws.Range("A40").Value = 45888
$ws.Range("A40").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B40").Value = 6411.3701171875
$ws.Range("C40").Value = 6456.47998046875
$ws.Range("D40").Value = 6400.22021484375
$ws.Range("E40").Value = 6446.240234375
$ws.Range("F40").Value = 4377440000
$ws.Range("G40").Value = -0.0058581031187567
$ws.Range("H40").Value = -1

$headlines = 'poriferous llc achieves coveted eu mdr regulatory certification for porous polyethylene implant warehouse data reveals how shipper managed inventory through tariff nvidia earnings face high bar on datacenter demand china deal newsmax to pay million to settle dominion suit over election fraud claim bmws next big suv is gunning for the top dog but will it go electric jules witcover political reporter and columnist dy at samsung overtakes apple u smartphone market update hint at major shift brent wti edge higher a zelensky arrives in washington for talk with trump top trending coin this week oz presale token gain early momentum best altcoins to buy for september remittix top list over ada pepe and bonk cold wallet offer referral bonus bonk burn t supply ripple eye crypto project set to explode top crypto executive seen accumulating million of token of new meme coin challenging shiba inu shib in walmart is selling a queen mattress topper for only thats soft and cool to the touch best electronic signature apps in judge denies bail to north andover offduty officer shot in home defense file another appeal whale bet big on this new payment token to replace xrp and stellar by carvanas breakout fuel interest in consumer ecommerce etf why are colorado rabbit growing tentacle and horn maluma sign with wme in all area the sky is the limit dow jones today stock tick lower a major index hover near record high after two consecutive week of gain investopedia fake labubu doll pose threat to kid safety watchdog say cbs news court of federal claim side with small business in landmark taa v baa protest employer guidance for new tax deduction for tip and overtime pay terra balcanica intersects gold mineralization over m and start drill program at brezani target in bosnia ap business summarybrief at pm edt conservative network newsmax agrees to pay m in defamation case over bogus election claim invesco qqq qqq share sold by lee financial co bank stock to keep an eye on august th terawulf stock jump more than a google boost stake in data center operator top bse gainer to watch in uncy investor have opportunity to lead unicycive therapeutic inc security fraud lawsuit with the schall law firm cornerstone select advisor llc reduces stock position in pfizer inc pfe southeast toyota finance ranked highest in overall dealer satisfaction by jd power forget ai robotic process automation is healthcare hottest technology warehouse teamster at unfi win lucrative first contract promising retail stock to follow now august th alphabet inc goog share sold by wellington management group llp promising lithium stock to keep an eye on august th insurance stock to watch today august th hydrogen stock to consider august th pwc kunai acquisition to bolster tech enablement for bank cfo the bull case for nvidia just got stronger according to these analyst this curved samsung screen is at one of it lowest price today soho house to be taken private in billion deal the new york time expect fiery outcome from fed jackson hole conference this week financial advisor network inc take position in cocacola company the ko what analyst rating have to say about the campbell mesa laboratory executes previouslyannounced strategic financing plan ibm option trading a deep dive into market sentiment forecasting the future analyst projection for cloudflare a closer look at walmarts option market dynamic flex meal deal return for faculty staff and graduate and professional student paypal holdingss option frenzy what you need to know a glimpse into the expert outlook on jm smucker through analyst assessing microsoft insight from financial analyst this ethereum trader turned into million in month the analyst verdict fedex in the eye of expert ai take apple edge foxconns pivot could spark etf shift where prudential financial stand with analyst mesa laboratory executes previouslyannounced strategic financing plan wired roundup why gpt flopped wired is the market bullish or bearish on royal gold looking into willis tower watson recent short interest a look into pge inc price over earnings how is the market feeling about lumentum holding ozempic creator novo nordisk will sell weight loss drug for cheaper through this new partnership novo nordisk offer ozempic at per month to eligible u cashpaying customer reuters texas may have to shut down data center to protect it energy grid can blockdag reach analyst review of the m presale in dogecoin news today dogecoin price decline but technical indicator suggest potential rebound applied digitals stock apld jump on plan for billion ai campus tipranks worker are job hugging in a stagnant labor market but growing resentment mean they could bail a soon a the next great resignation come yahoocom texas attorney general accuses meta characterai of misleading kid with mental health claim techcrunch worker are job hugging in a stagnant labor market but growing resentment mean they could bail a soon a the next great resignation come fortune conservative network newsmax agrees to pay m in defamation case over bogus election claim skechers named title sponsor of world champion cup shriners childrens announced a official charity partner feather sound country club in clearwater florida to host tournament whats going on with propanc share monday mastercard and worldpay team to streamline uae money movement mcdonalds corporation mcd share sold by community bank trust waco texas granite harbor advisor inc boost holding in mcdonalds corporation mcd walmart ha everything you need to head back to school from backpack to laptop madd join force with law enforcement to prevent impaired driving ahead of labor day with saturation saturday initiative this weekend ibd live qa stock list for monday aug spotlight on fed what powell bowman and waller will say about rate dogecoin rally spark buzz but trader say this sub coin ha more upside fairscale capital llc purchase share of crowdstrike crwd lee financial co sell share of danaher corporation dhr berger montague pc investigates security claim against kindercare learning company inc nyse klc amazon is selling an bluetooth speaker for and it perfect for any setting inventhelp inventor develops new training aid for basketball player chk top data science recruiter in find ai analytics and ml expert gsg advisor llc purchase share of invesco qqq qqq gsg advisor llc ha stock holding in mastercard incorporated ma weekly bankruptcy alert august for the week ending august jackson hole loom a powell speech test september ratecut hope newsmax to pay million to settle dominion defamation lawsuit intel stock drop on report of equity stake via the chip act barrons liquid edge ai platform leap expands support to laptop with bestinclass performance on amd ryzentm and ryzen aitm processor ap business summarybrief at pm edt ap business summarybrief at pm edt texas democrat are back what now froma harrop are parent who dont vaccinate their child guilty of abuse handwritten proof of holocaust theft should compel congress to act ap news summary at pm edt ap news summary at pm edt state and dc sue doj to stop immigration requirement on victim fund trump zelenskyy meeting not end of the road for u support in securing a peace deal editorial roundup united state eric adam roll out nypd expansion say trump dc crackdown not a factor ap trending summarybrief at pm edt texas lawmaker in illinois heading back home following twoweek walkout over remap fight nfl team refuse to back down after facing backlash for hiring male cheerleader stephen miller say trump dc crackdown will also target graffiti ap news summary at pm edt conservative network newsmax agrees to pay m in defamation case over bogus election claim mark hamill say his wife clever thought stopped him from fleeing u after trump won trump zelensky put on the smile at oval office vancouver mayor warns delay on i bridge replacement cost m per day trump zelensky put on the smile at oval office hamas accepts an arab ceasefire proposal on gaza a palestinian death toll pass texas redistricting fight shift a democrat return to austin with a new strategy mississippi becomes fourth state to send national guard troop to dc in expanding federal crackdown mississippi becomes fourth state to send national guard troop to dc in expanding federal crackdown mississippi becomes fourth state to send national guard troop to dc in expanding federal crackdown mississippi becomes fourth state to send national guard troop to dc in expanding federal crackdown mississippi becomes fourth state to send national guard troop to dc in expanding federal crackdown mississippi becomes fourth state to send national guard troop to dc in expanding federal crackdown mississippi becomes fourth state to send national guard troop to dc in expanding federal crackdown mississippi becomes fourth state to send national guard troop to dc in expanding federal crackdown mississippi becomes fourth state to send national guard troop to dc in expanding federal crackdown mississippi becomes fourth state to send national guard troop to dc in expanding federal crackdown on putin advice trump launch assault on mailin ballot and voting machine eo seek to expand access to alternative asset investment in retirement plan doj demand sensitive illinois voter registration data mike mcgrew an alternative mindset for patriotic independent commentary ha brian stelter doesnt get that he only making msnbcs rebranding even funnier in his excited thread former attorney general william barr say he didnt see trump implicated in epstein file former attorney general william barr say he didnt see trump implicated in epstein file gavin newsom left fox news speechless with his unexpectedly hilarious woke comeback these state leader might be the last line of defense against trump abuse of power subsidence isnt just an environmental crisis and it can be slowed the latest zelenskyy signal openness to threeway meeting with trump and putin doj wont appeal judge order on concealed carry ban in post office the latest zelenskyy signal openness to threeway meeting with trump and putin letter political treachery versus personal charity what to know about redistricting fight a texas democrat return and california start work what to know about redistricting fight a texas democrat return and california start work what to know about redistricting fight a texas democrat return and california start work what to know about redistricting fight a texas democrat return and california start work what to know about redistricting fight a texas democrat return and california start work what to know about redistricting fight a texas democrat return and california start work what to know about redistricting fight a texas democrat return and california start work what to know about redistricting fight a texas democrat return and california start work what to know about redistricting fight a texas democrat return and california start work bnpl provider ask for extension to nydfs request for information deadline wyoming teen educate youth about human trafficking on social medium congress must correct an injustice facing combatinjured veteran how blessing and disaster shape alabama why is the c galaxy such an iconic aircraft texas attorney general accuses meta characterai of misleading kid with mental health claim texas ag accuses meta characterai of misleading kid with mental health claim european leader join trump zelenskyy at white house conservative network newsmax agrees to pay m in defamation case over bogus election claim what to know about redistricting fight a texas democrat return and california start work what to know about redistricting fight a texas democrat return and california start work twoweek standoff end quorum reached in texas house mississippi becomes fourth state to send national guard troop to dc in expanding federal crackdown letter trumpputin photo op accomplished nothing judge weighs detainee legal right at alligator alcatraz in florida everglades judge weighs detainee legal right at alligator alcatraz in florida everglades roblox is locking down sexual content and access to adult location after lawsuit trump say the u would back european security guarantee for ukraine to help end the war with russia good in a crisis dr nirav shah could be a promising candidate for anything steve collins wear a suit bro zelenskyy greeted by banner begging him to dress for the occassion did he hulk hogan couldnt stand gavin newsom exwife say about govs online trolling rightwing reporter brian glenn praise zelenskys attire after shaming him in february texas democrat return ending twoweek walkout over redistricting trump say hell talk to putin after meeting with zelenskyy and european leader at the white house reporter brian glenn brag about dc being safe after walking street with girlfriend marjorie taylor greene zelenskyy trump express hope for trilateral talk with putin to bring end to russiaukraine war oval office meeting again put rocky tense trumpzelenskyy relationship to the test armstrong williams how george washington would end russia war on ukraine staff commentary la habra man charged in fatal whittier shooting of driver texas democrat head home after week in illinois mpo holding first public planning session tuesday at mountain line photo show aftermath of a deadly russian attack in ukraine and the funeral of a fallen soldier now thats a selfown eric swalwell grocery shopping go viral for all the wrong hilarious reason texas democrat return to capitol for republican u house map overhaul a california counter texas democrat return to capitol for republican u house map overhaul a california counter letter time for democrat to play to win volodymyr zelensky wear suit this time for white house meeting with trump tarrant county commissioner again propose cutting voting site ahead of an election this is exactly what technology should be doing man shot injured in parking lot of huntington beach store letter no place for threat violence in politics trump aim to set putinzelenskiy meeting after monday summit obbba creates a cohesive framework for u hydrogen growth top trump lackey demand dc resident kiss the ground for occupation republican look to make a uturn on electric vehicle for the postal service republican look to make a uturn on electric vehicle for the postal service'
$ws.Range("I40").Value = $headlines
$ws.Range("J40").Value = $headlines

$ws.Range("K40").Value = 0.13
$ws.Range("L40").Value = 0.106
$ws.Range("M40").Value = 0.764
$ws.Range("N40").Value = 0.9955000000000001
$ws.Range("O40").Value = 94
$ws.Range("P40").Value = 78
$ws.Range("Q40").Value = 0
$ws.Range("R40").Value = 140
$ws.Range("S40").Value = 62
$ws.Range("T40").Value = 203
$ws.Range("U40").Value = 113
$ws.Range("V40").Value = 49
$ws.Range("W40").Value = 35
$ws.Range("X40").Value = 64

